# Insert a new data row at row 102 (shifts existing rows 102-213 down to 103-214)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with its full record.
# (Column layout mirrors every other data row in the sheet.)
$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44789
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100104
$ws.Range("H102").Value = "Frutos de pepita"
$ws.Range("I102").Value = 100104003
$ws.Range("J102").Value = "Membrillo"
$ws.Range("K102").Value = "Champion"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 100
$ws.Range("N102").Value = 10000
$ws.Range("O102").Value = 10000
$ws.Range("P102").Value = 10000
$ws.Range("Q102").Value = "$/bandeja 18 kilos granel"
$ws.Range("R102").Value = "Región de O'Higgins"
$ws.Range("S102").Value = 556
$ws.Range("T102").Value = 18

# Match the date-cell style/number format used by the rest of column D.
$ws.Range("D102").NumberFormat = $ws.Range("D103").NumberFormat
